# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 17 de Octubre de 2020 a las 07:13"

# Row 60 (Uzbekistan) - updated case numbers
$ws.Range("B60").Value = 62684
$ws.Range("C60").Value = 96
$ws.Range("D60").Value = 59756
$ws.Range("E60").Value = 2408

# Rows 146/147 - Tailandia overtakes Guyana in ranking (sorted desc by Casos totales)
$ws.Range("A146").Value = "Tailandia"
$ws.Range("B146").Value = 3679
$ws.Range("C146").Value = 10
$ws.Range("D146").Value = 3478
$ws.Range("E146").Value = 142
$ws.Range("H146").Value = 59

$ws.Range("A147").Value = "Guyana"
$ws.Range("B147").Value = 3672
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 2590
$ws.Range("E147").Value = 975
$ws.Range("H147").Value = 107

# Rows 216/217 - Islas Malvinas and Montserrat swap ranking order
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0

$ws.Range("A217").Value = "Montserrat"
$ws.Range("D217").Value = 12
$ws.Range("H217").Value = 1
